# Edit: resolve bug with empty footnotes
# - footnote 23: drop stray trailing "aa"
# - footnote 42: fill in previously-empty note text
# - footnote 63: remove entirely (was an empty/orphan note), along with
#   the trailing space run that preceded its reference in the body

$d = $word.ActiveDocument

# 1) Footnote id=23 (collection index 3): remove stray "aa" suffix.
$fn23 = $d.Footnotes.Item(3)
$fn23.Range.Text = $fn23.Range.Text.Replace("གྱིས། སྣར་ཐང་། པེ་ཅིན།aa", "གྱིས། སྣར་ཐང་། པེ་ཅིན།")

# 2) Footnote id=42 (collection index 22): was just a lone danda "।",
#    fill in the real note text.
$fn42 = $d.Footnotes.Item(22)
$fn42.Range.Text = "ཉིད། ཞེས་པར་མ་གཞན་ནང་མེད།"

# 3) Footnote id=63 (collection index 43, the last footnote): delete the
#    whole (empty/orphan) note - removes both the <w:footnote> body and
#    its <w:footnoteReference> run in the body.
$fn63 = $d.Footnotes.Item(43)
$fn63.Delete()

# 4) The footnote reference for id=63 was preceded by a run containing a
#    single space; with the reference gone that trailing space run is an
#    orphan too, so trim it from the end of the (single) body paragraph.
$c = $d.Content
$trailingSpace = $d.Range($c.End - 2, $c.End - 1)
$trailingSpace.Delete()
